$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs | Ccl21b | Cxcr3 | FAPs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.134289
$ws.Range("H2").Value = 0.402867
$ws.Range("I2").Value = 0.3678949098679525
$ws.Range("J2").Value = 0.3678949098679525
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05194133333333333
$ws.Range("N2").Value = 0.155824
$ws.Range("O2").Value = 0.03788844568234288
$ws.Range("P2").Value = 0.03788844568234288
$ws.Range("Q2").Value = 0.006975149712000001
$ws.Range("R2").Value = 0.06277634740800001
$ws.Range("S2").Value = 0.01393896630934235
$ws.Range("T2").Value = 0.01393896630934235

# Row 3: FAPs | Ccl21b | Cxcr3 | M2
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.134289
$ws.Range("H3").Value = 0.402867
$ws.Range("I3").Value = 0.3678949098679525
$ws.Range("J3").Value = 0.3678949098679525
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.318960333333333
$ws.Range("N3").Value = 3.956881
$ws.Range("O3").Value = 0.962111554317657
$ws.Range("P3").Value = 0.9621115543176572
$ws.Range("Q3").Value = 0.177121864203
$ws.Range("R3").Value = 1.594096777827
$ws.Range("S3").Value = 0.3539559435586101
$ws.Range("T3").Value = 0.3539559435586102

# Row 4: sCs | Ccl21b | Cxcr3 | FAPs
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Ccl21b"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.230731
$ws.Range("H4").Value = 0.6921929999999999
$ws.Range("I4").Value = 0.6321050901320475
$ws.Range("J4").Value = 0.6321050901320475
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05194133333333333
$ws.Range("N4").Value = 0.155824
$ws.Range("O4").Value = 0.03788844568234288
$ws.Range("P4").Value = 0.03788844568234288
$ws.Range("Q4").Value = 0.01198447578133333
$ws.Range("R4").Value = 0.107860282032
$ws.Range("S4").Value = 0.02394947937300054
$ws.Range("T4").Value = 0.02394947937300054

# Row 5: sCs | Ccl21b | Cxcr3 | M2
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Ccl21b"
$ws.Range("C5").Value = "Cxcr3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.230731
$ws.Range("H5").Value = 0.6921929999999999
$ws.Range("I5").Value = 0.6321050901320475
$ws.Range("J5").Value = 0.6321050901320475
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.318960333333333
$ws.Range("N5").Value = 3.956881
$ws.Range("O5").Value = 0.962111554317657
$ws.Range("P5").Value = 0.9621115543176572
$ws.Range("Q5").Value = 0.3043250366703333
$ws.Range("R5").Value = 2.738925330033
$ws.Range("S5").Value = 0.6081556107590469
$ws.Range("T5").Value = 0.608155610759047
